$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Team" -> "Payment By"
$ws.Range("H1").Value = "Payment By"

# Row 2 (Expense #8)
$ws.Range("A2").Value = "#00-8"
$ws.Range("C2").Value = "Anuj Pal"
$ws.Range("D2").Value = "My first expense"
$ws.Range("E2").Value = "Travel Expense"
$ws.Range("F2").Value = "Attached"
$ws.Range("G2").Value = 2250
$ws.Range("H2").Value = "Hard Cash"

# Row 3 (Expense #9)
$ws.Range("A3").Value = "#00-9"
$ws.Range("C3").Value = "Anuj Pal"
$ws.Range("D3").Value = "My first expense"
$ws.Range("E3").Value = "Travel Expense"
$ws.Range("F3").Value = "Attached"
$ws.Range("G3").Value = 12250
$ws.Range("H3").Value = "Hard Cash"

# Row 4 (Expense #10)
$ws.Range("A4").Value = "#00-10"
$ws.Range("C4").Value = "Anuj Pal"
$ws.Range("D4").Value = "fdgds"
$ws.Range("E4").Value = "Food Expense"
$ws.Range("F4").Value = "Attached"
$ws.Range("G4").Value = 4355
$ws.Range("H4").Value = "Debit Card"
